$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Operation 1: the "阴晴圆缺" paragraph (#7) loses its paragraph-mark
# formatting (the <w:pPr> carrying rFonts hint="eastAsia"), while its run
# content stays the same. To redefine a paragraph mark via InsertXML we
# must select the whole paragraph plus the first character of the
# following paragraph ("2022.07.02"), then restore that borrowed
# character (with its original formatting) as part of the replacement.
$p7 = $d.Paragraphs(7)
$r1 = $d.Range($p7.Range.Start, $p7.Range.End + 1)
$xml1 = "<w:p $wNs><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>阴晴圆缺</w:t></w:r></w:p><w:p $wNs><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t>022.07.02</w:t></w:r></w:p>"
$r1.InsertXML($xml1) | Out-Null

# --- Operation 2: the final paragraph ("阴晴不定") is split in two. The
# first half keeps the plain text "阴晴不定" with no special paragraph
# mark formatting. A new paragraph is appended after it which carries the
# original paragraph-mark formatting (rFonts hint="eastAsia") and holds
# the new text "12345" (entered as "1" then "2345").
$last = $d.Paragraphs($d.Paragraphs.Count)
$r2 = $d.Range($last.Range.Start, $last.Range.End)
$xml2 = "<w:p $wNs><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>阴晴不定</w:t></w:r></w:p><w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>1</w:t></w:r><w:r><w:t>2345</w:t></w:r></w:p>"
$r2.InsertXML($xml2)
